$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The first three summary cells (row 1-3) all become "0M".
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# The per-iteration values that used to live tab-separated inside a single
# cell (near the end of the table) now get their own one-value-per-row
# table rows, inserted right after row 3.
$newRowTexts = @(
    "101",
    "0.00003",
    "0.00005",
    "0.00004",
    "0.00000",
    "0.00004",
    "0.00004",
    "0.00004",
    "0.00373",
    "100.0"
)

$afterRow = 3
foreach ($txt in $newRowTexts) {
    $afterRow = $afterRow + 1
    $newRow = $t.Rows.Add($t.Rows.Item($afterRow))
    $newRow.Cells.Item(1).Range.Text = $txt
}

# The two rows that used to hold the tab-separated per-iteration values are
# collapsed down to a single summary value each, and the trailing empty row
# picks up the summary value that used to be in row 3.
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.98"
$t.Cell($rowCount - 1, 1).Range.Text = "0"
$t.Cell($rowCount, 1).Range.Text = "15"
